$wb = $excel.ActiveWorkbook

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = 0.006565603454886704
$wsLinear.Range("B3").Value = -0.06374982851663631
$wsLinear.Range("B4").Value = 0.1717260184579955
$wsLinear.Range("B5").Value = "[1.0, 0.16198153230227647, 0.04495023656944202, 0.007773353046522208, 0.010398346581749868, 0.013934927903101108, 0.07430826103259126, 0.15850968380101207, 0.06013840848908155, -0.0034959827032100747, -0.04499710160011272, -0.049903346542983416, -0.05333345675419163, 0.0176190420107261, 0.07222439465907735, 0.038993783645925616, -0.023107413638629896, -0.05390159183813175, -0.008204865040037713, -0.022093481660576602]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = -0.003795207957694611
$wsNonLinear.Range("B5").Value = -0.1280207906572329
$wsNonLinear.Range("B6").Value = 0.188911286317794
$wsNonLinear.Range("B7").Value = 0.004891879908572727
$wsNonLinear.Range("B8").Value = -0.02815402019901062
$wsNonLinear.Range("B9").Value = 0.1578207405216499
$wsNonLinear.Range("B10").Value = "[0.9999999999999999, 0.16308746550477266, 0.046932344942685965, 0.007281137601957118, 0.010912954082495922, 0.014522615722712405, 0.07368306220441415, 0.1593403889608029, 0.06017055512045514, -0.003846349573732074, -0.0466000327892286, -0.051921877701179765, -0.054188352563159126, 0.01685643792364399, 0.07281896278065549, 0.039228977153699114, -0.023990719274032225, -0.054968107644284364, -0.009682270857879519, -0.0227657542104622]"
